$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the specific project code with a generic placeholder
$ws.Range("C2").Value = "[DL-MAA20XX-YY]"

$ws.Range("C2").Select()
